$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-10"

# Update the header label cell (A1 shared string "2022 (through 06-09)")
$ws.Range("I1").Value = "2022 (through 06-10)"

# Update June value (row 7) and Total value (row 14) in column I
$ws.Range("I7").Value = 35
$ws.Range("I14").Value = 698
